$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where an "N" (week 12 presence) boolean TRUE value must be added.
$rows = @(10, 11, 18, 21, 25, 28, 36, 37, 41)

foreach ($r in $rows) {
    $ws.Range("N$r").Value = $true
}

# Update the view: remove the frozen/scrolled top-left cell and move the
# active selection to O6 (matches the saved selection state in the file).
$ws.Range("O6").Select()

$wb.Save()
